# Bug fix - broken immutable index when formatting, format index now
# always creates new index.
#
# The quarterly performance table on slide 1 gets its "Group1" header
# (first of the two merged two-column group headers) renamed to
# "Group2", all of the percentage figures in the body recalculated, and
# the "2017"/"2017" year labels on the last two data rows corrected to
# "2018".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$sh  = $s.Shapes.Item(2)
$tbl = $sh.Table

# Header row: first "Group1"/"Group1" merged pair -> "Group2"
$tbl.Cell(1,3).Shape.TextFrame.TextRange.Text = "Group2"
$tbl.Cell(1,4).Shape.TextFrame.TextRange.Text = "Group2"

# Mar 31, 2017 row
$tbl.Cell(3,3).Shape.TextFrame.TextRange.Text = "0.45%"
$tbl.Cell(3,4).Shape.TextFrame.TextRange.Text = "0.04%"
$tbl.Cell(3,5).Shape.TextFrame.TextRange.Text = "0.96%"
$tbl.Cell(3,6).Shape.TextFrame.TextRange.Text = "0.93%"
$tbl.Cell(3,7).Shape.TextFrame.TextRange.Text = "2.37%"

# Jun 30, 2017 row
$tbl.Cell(4,3).Shape.TextFrame.TextRange.Text = "0.08%"
$tbl.Cell(4,4).Shape.TextFrame.TextRange.Text = "0.02%"
$tbl.Cell(4,5).Shape.TextFrame.TextRange.Text = "0.12%"
$tbl.Cell(4,6).Shape.TextFrame.TextRange.Text = "0.71%"
$tbl.Cell(4,7).Shape.TextFrame.TextRange.Text = "0.94%"

# Sep 30, 2017 row
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Text = "0.12%"
$tbl.Cell(5,4).Shape.TextFrame.TextRange.Text = "0.34%"
$tbl.Cell(5,5).Shape.TextFrame.TextRange.Text = "0.78%"
$tbl.Cell(5,6).Shape.TextFrame.TextRange.Text = "0.24%"
$tbl.Cell(5,7).Shape.TextFrame.TextRange.Text = "1.48%"

# Dec 31, 2017 row
$tbl.Cell(6,3).Shape.TextFrame.TextRange.Text = "0.66%"
$tbl.Cell(6,4).Shape.TextFrame.TextRange.Text = "0.51%"
$tbl.Cell(6,5).Shape.TextFrame.TextRange.Text = "0.37%"
$tbl.Cell(6,6).Shape.TextFrame.TextRange.Text = "0.27%"
$tbl.Cell(6,7).Shape.TextFrame.TextRange.Text = "1.82%"

# Mar 31, 2018 row (Year label corrected 2017 -> 2018)
$tbl.Cell(7,1).Shape.TextFrame.TextRange.Text = "2018"
$tbl.Cell(7,3).Shape.TextFrame.TextRange.Text = "0.76%"
$tbl.Cell(7,4).Shape.TextFrame.TextRange.Text = "0.93%"
$tbl.Cell(7,5).Shape.TextFrame.TextRange.Text = "0.17%"
$tbl.Cell(7,6).Shape.TextFrame.TextRange.Text = "0.95%"
$tbl.Cell(7,7).Shape.TextFrame.TextRange.Text = "2.80%"

# Jun 30, 2018 row (Year label corrected 2017 -> 2018)
$tbl.Cell(8,1).Shape.TextFrame.TextRange.Text = "2018"
$tbl.Cell(8,3).Shape.TextFrame.TextRange.Text = "0.01%"
$tbl.Cell(8,4).Shape.TextFrame.TextRange.Text = "0.18%"
$tbl.Cell(8,5).Shape.TextFrame.TextRange.Text = "0.66%"
$tbl.Cell(8,6).Shape.TextFrame.TextRange.Text = "0.84%"
$tbl.Cell(8,7).Shape.TextFrame.TextRange.Text = "1.70%"

# Total row
$tbl.Cell(9,3).Shape.TextFrame.TextRange.Text = "2.08%"
$tbl.Cell(9,4).Shape.TextFrame.TextRange.Text = "2.02%"
$tbl.Cell(9,5).Shape.TextFrame.TextRange.Text = "3.05%"
$tbl.Cell(9,6).Shape.TextFrame.TextRange.Text = "3.95%"
$tbl.Cell(9,7).Shape.TextFrame.TextRange.Text = "11.09%"
